$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 0. Seed the shared-string table in the exact order the target file
#    uses: RADIANS, DEGREES, POWER, MOD, SQRT (indices 21..25).
# ---------------------------------------------------------------------
$ws.Range("A46").Value = "RADIANS"
$ws.Range("A45").Value = "DEGREES"
$ws.Range("A47").Value = "POWER"
$ws.Range("A48").Value = "MOD"
$ws.Range("A49").Value = "SQRT"

# ---------------------------------------------------------------------
# 1. Row 34 (ATAN2 example) - new inputs + shared formula group J34:M34
# ---------------------------------------------------------------------
$ws.Range("C34").Value = -1
$ws.Range("D34").Value = 2.5
$ws.Range("E34").Value = -10
$ws.Range("J34:M34").Formula = "=ATAN2(C34,D34)"

# ---------------------------------------------------------------------
# 2. Row 45 - DEGREES example
# ---------------------------------------------------------------------
$ws.Range("B45").Formula = "=PI()"
$ws.Range("C45").Formula = "=-PI()/4"
$ws.Range("D45").Value = 120
$ws.Range("F45").Value = "A"
$ws.Range("I45").Formula = "=DEGREES(B45)"
$ws.Range("J45").Formula = "=DEGREES(C45)"
$ws.Range("K45").Formula = "=DEGREES(D45)"
$ws.Range("L45").Formula = "=DEGREES(E45)"
$ws.Range("M45").Formula = "=DEGREES(F45)"
$ws.Range("O45").FormulaArray = "=SUM(DEGREES(B45:E45))"

# ---------------------------------------------------------------------
# 3. Row 46 - RADIANS example (shared formula group J46:M46)
# ---------------------------------------------------------------------
$ws.Range("B46").Value = 120
$ws.Range("C46").Value = -45
$ws.Range("D46").Value = -720
$ws.Range("F46").Value = "A"
$ws.Range("I46").Formula = "=RADIANS(B46)"
$ws.Range("J46:M46").Formula = "=RADIANS(C46)"
$ws.Range("O46").FormulaArray = "=SUM(RADIANS(B46:E46))"

# ---------------------------------------------------------------------
# 4. Row 47 - POWER example (replaces the old single styled F47 cell)
# ---------------------------------------------------------------------
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = 3.5
$ws.Range("D47").Value = -0.5
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = "A"
$ws.Range("F47").Style = "Normal"
$ws.Range("I47").Formula = "=POWER(B47,C47)"
$ws.Range("J47").Formula = "=POWER(C47,D47)"
$ws.Range("K47").Formula = "=POWER(D47,E47)"
$ws.Range("L47").Formula = "=POWER(E47,F47)"
$ws.Range("M47").Formula = "=POWER(F47,G47)"
$ws.Range("O47").FormulaArray = "=SUM(POWER(B47:D47,C47:E47))"

# ---------------------------------------------------------------------
# 5. Row 48 - MOD example
# ---------------------------------------------------------------------
$ws.Range("B48").Value = 31.5
$ws.Range("C48").Value = 10
$ws.Range("D48").Value = 6.7
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = "A"
$ws.Range("I48").Formula = "=MOD(B48,C48)"
$ws.Range("J48").Formula = "=MOD(C48,D48)"
$ws.Range("K48").Formula = "=MOD(D48,E48)"
$ws.Range("L48").Formula = "=MOD(E48,F48)"
$ws.Range("M48").Formula = "=MOD(F48,G48)"
$ws.Range("O48").FormulaArray = "=SUM(MOD(B48:D48,C48:E48))"

# ---------------------------------------------------------------------
# 6. Row 49 - SQRT example (shared formula group J49:M49)
# ---------------------------------------------------------------------
$ws.Range("B49").Value = 2
$ws.Range("C49").Value = 3.5
$ws.Range("D49").Value = 9
$ws.Range("E49").Value = -0.5
$ws.Range("F49").Value = "A"
$ws.Range("I49").Formula = "=SQRT(B49)"
$ws.Range("J49:M49").Formula = "=SQRT(C49)"
$ws.Range("O49").FormulaArray = "=SUM(SQRT(B49:D49))"

# ---------------------------------------------------------------------
# 7. Rows 52-60 - the wrap-text styled filler cells (moved down from the
#    old row 47).
# ---------------------------------------------------------------------
foreach ($r in 52..58) {
    $ws.Range("F$r").WrapText = $true
    $ws.Range("H$r").WrapText = $true
}
$ws.Range("E55").WrapText = $true
$ws.Range("H59").WrapText = $true
$ws.Range("H60").WrapText = $true

# ---------------------------------------------------------------------
# 8. Column widths for F:G (best-fit, matches new wider "A"/formula text)
# ---------------------------------------------------------------------
$ws.Columns("F:G").AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 9. View state - scroll down to the new block and select it, matching
#    the author's saved cursor position.
# ---------------------------------------------------------------------
$ws.Rows("51:66").Select()
